{"js": "// Change \"... redistributed to other pits in clockwise direction.\" to\n// \"... redistributed to other pits in anti-clockwise direction.\" by\n// inserting \"anti-\" immediately before the (unique) word \"clockwise\".\n//\n// Using a targeted search + insertText keeps the existing run formatting\n// (green highlight, en-GB language) intact, because Word.js inherits the\n// surrounding run's character formatting for inserted text.\nconst searchResults = context.document.body.search(\"clockwise direction.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find \"clockwise direction.\" in the document body.');\n}\n\n// Insert \"anti-\" right before every matched occurrence (i.e. right before\n// \"clockwise\"). In this document the phrase is unique, so this updates the\n// single \"5.2 ...\" sentence, but looping keeps the script correct even if\n// the phrase were repeated.\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"anti-\", Word.InsertLocation.before);\n}\nawait context.sync();\n", "ps1": "# Change \"... redistributed to other pits in clockwise direction.\" to\n# \"... redistributed to other pits in anti-clockwise direction.\"\n#\n# A plain Find & Replace on the unique phrase keeps the existing run\n# formatting (green highlight, en-GB language) untouched, since only the\n# text inside the matched range is substituted.\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Text = \"clockwise direction.\"\n$rng.Find.Replacement.Text = \"anti-clockwise direction.\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 1\n\n# wdReplaceOne = 1 -> replace just the single match found\n$found = $rng.Find.Execute($rng.Find.Text, $rng.Find.MatchCase, $rng.Find.MatchWholeWord, `\n    $false, $false, $false, $rng.Find.Forward, $rng.Find.Wrap, $false, `\n    $rng.Find.Replacement.Text, 1)\n\nif (-not $found) {\n    Write-Output \"WARNING: 'clockwise direction.' was not found; no replacement made.\"\n}\n"}
